$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns L:N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the bold/bordered header formatting already used by A1:K1
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns L:N for rows 2-7
$data = @(
    @(91.17116223450478, 202460, 321.8759936406995),
    @(91.74154684374953, 61093, 303.9452736318408),
    @(90.6205976008235, 151336, 141.7003745318352),
    @(88.85439072913162, 52310, 166.0634920634921),
    @(18.28712748796549, 1935, 14.33333333333333),
    @(28.57643247462115, 273, 16.05882352941176)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $data[$i][0]
    $ws.Cells.Item($row, 13).Value = $data[$i][1]
    $ws.Cells.Item($row, 14).Value = $data[$i][2]
}
